# Apply updated cryptocurrency price/volume figures scraped on
# Wed Sep 13 08:41:33 UTC 2023, matching the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a pure-text value (matching the workbook's
    # inlineStr cells) even when the string looks like a number, then
    # drop the temporary text format so the cell style stays default.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range('D2').Value = '26.022.32'
$ws.Range('E2').Value = '  +0.51%  '

# Row 3
$ws.Range('D3').Value = '1.596.47'
$ws.Range('E3').Value = '  +0.76%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
Set-TextValue $ws.Range('D5') '211.89'
$ws.Range('E5').Value = '  +0.80%  '

# Row 6
$ws.Range('E6').Value = '  -0.08%  '

# Row 7
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('E8').Value = '  +0.35%  '

# Row 9
$ws.Range('E9').Value = '  +0.10%  '

# Row 10
Set-TextValue $ws.Range('D10') '18.25'
$ws.Range('E10').Value = '  +0.93%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0809'
$ws.Range('E11').Value = '  +2.35%  '

# Row 12
$ws.Range('D12').Value = '1.819.32'
$ws.Range('E12').Value = '  +0.81%  '

# Row 13
$ws.Range('D13').Value = '1.588.46'
$ws.Range('E13').Value = '  +0.20%  '

# Row 14
$ws.Range('E14').Value = '  -0.37%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.516'
$ws.Range('E15').Value = '  +1.91%  '

# Row 16
$ws.Range('D16').Value = '26.014.50'
$ws.Range('E16').Value = '  +0.49%  '

# Row 17
Set-TextValue $ws.Range('D17') '60.78'
$ws.Range('E17').Value = '  +1.32%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.56%  '

# Row 19
$ws.Range('E19').Value = '  -0.09%  '

# Row 20
Set-TextValue $ws.Range('D20') '203.38'
$ws.Range('E20').Value = '  +5.36%  '

# Row 21
$ws.Range('E21').Value = '  +1.62%  '

# Row 22
Set-TextValue $ws.Range('D22') '9.27'
$ws.Range('E22').Value = '  -0.99%  '

# Row 23
Set-TextValue $ws.Range('D23') '6.05'
$ws.Range('E23').Value = '  +1.86%  '

# Row 24
Set-TextValue $ws.Range('D24') '1.93'
$ws.Range('E24').Value = '  +13.20%  '

# Row 25
Set-TextValue $ws.Range('D25') '143.82'
$ws.Range('E25').Value = '  +1.69%  '

# Row 26
$ws.Range('E26').Value = '  -0.08%  '

# Row 27
$ws.Range('E27').Value = '  -7.54%  '

# Row 28
$ws.Range('E28').Value = '  +0.61%  '

# Row 29
Set-TextValue $ws.Range('D29') '6.53'
$ws.Range('E29').Value = '  +1.37%  '

# Row 30
$ws.Range('E30').Value = '  +0.68%  '

# Row 31
$ws.Range('E31').Value = '  +1.03%  '

# Row 32
$ws.Range('E32').Value = '  +0.36%  '

# Row 33
Set-TextValue $ws.Range('D33') '2.91'
$ws.Range('E33').Value = '  -3.71%  '

# Row 34
$ws.Range('E34').Value = '  -0.58%  '

# Row 35
$ws.Range('E35').Value = '  -0.59%  '

# Row 36
$ws.Range('D36').Value = '1.130.22'
$ws.Range('E36').Value = '  +3.10%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.0164'
$ws.Range('E37').Value = '  +8.60%  '

# Row 38
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D38') '1.00'
$ws.Range('E38').Value = '  -0.10%  '

# Row 39
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D39') '0.798'
$ws.Range('E39').Value = '  +2.32%  '

# Row 40
$ws.Range('E40').Value = '  -0.86%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.493'
$ws.Range('E41').Value = '  -1.85%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.780'
$ws.Range('E42').Value = '  -2.09%  '

# Row 43
Set-TextValue $ws.Range('D43') '5.16'
$ws.Range('E43').Value = '  +0.56%  '

# Row 44
$ws.Range('D44').Value = '1.730.03'
$ws.Range('E44').Value = '  +0.63%  '

# Row 45
Set-TextValue $ws.Range('D45') '92.26'
$ws.Range('E45').Value = '  -1.29%  '

# Row 46
Set-TextValue $ws.Range('D46') '54.04'
$ws.Range('E46').Value = '  +1.70%  '

# Row 47
$ws.Range('E47').Value = '  -0.98%  '

# Row 48
$ws.Range('E48').Value = '  -0.45%  '

# Row 49
$ws.Range('B49').Value = 'USDD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws.Range('D49') '1.01'
$ws.Range('E49').Value = '  +0.45%  '

# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D50') '0.406'
$ws.Range('E50').Value = '  -0.27%  '

# Row 51
$ws.Range('D51').Value = '0.0₇0949'
$ws.Range('E51').Value = '  -15.22%  '
